# Update 2023 (column J) crime-count figures to reflect data through 2023-05-15
# across the Citywide Totals, By Neighborhood, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 2553
$ws.Range("J3").Value = 2618
$ws.Range("J4").Value = 601
$ws.Range("J5").Value = 200
$ws.Range("J6").Value = 3261
$ws.Range("J7").Value = 9233

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J2").Value = 35
$ws.Range("J7").Value = 106

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J3").Value = 110
$ws.Range("J7").Value = 309

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J4").Value = 7
$ws.Range("J7").Value = 127

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J4").Value = 23
$ws.Range("J6").Value = 101
$ws.Range("J7").Value = 334

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J6").Value = 16
$ws.Range("J7").Value = 69

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J2").Value = 72
$ws.Range("J6").Value = 84
$ws.Range("J7").Value = 238

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J5").Value = 23
$ws.Range("J7").Value = 275
$ws.Range("J8").Value = 577
$ws.Range("J11").Value = 130
$ws.Range("J19").Value = 298
$ws.Range("J21").Value = 14
$ws.Range("J23").Value = 96
$ws.Range("J29").Value = 537
$ws.Range("J31").Value = 69
$ws.Range("J33").Value = 376
$ws.Range("J36").Value = 137
$ws.Range("J37").Value = 309
$ws.Range("J42").Value = 359
$ws.Range("J43").Value = 84
$ws.Range("J44").Value = 75
$ws.Range("J48").Value = 91
$ws.Range("J53").Value = 88
$ws.Range("J57").Value = 43
$ws.Range("J63").Value = 46
$ws.Range("J65").Value = 238
$ws.Range("J67").Value = 334
$ws.Range("J72").Value = 34
$ws.Range("J73").Value = 80
$ws.Range("J78").Value = 128
$ws.Range("J79").Value = 281
$ws.Range("J85").Value = 430
$ws.Range("J86").Value = 55
$ws.Range("J88").Value = 93
$ws.Range("J90").Value = 101
$ws.Range("J91").Value = 104
$ws.Range("J92").Value = 29
$ws.Range("J94").Value = 77
$ws.Range("J96").Value = 106
$ws.Range("J97").Value = 57
$ws.Range("J99").Value = 127
$ws.Range("J101").Value = 9233

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 100
$ws.Range("J3").Value = 112
$ws.Range("J6").Value = 129
$ws.Range("J7").Value = 376

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 158
$ws.Range("J3").Value = 177
$ws.Range("J6").Value = 149
$ws.Range("J7").Value = 537

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J3").Value = 81
$ws.Range("J7").Value = 298

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J2").Value = 28
$ws.Range("J7").Value = 75

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J6").Value = 44
$ws.Range("J7").Value = 91

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J3").Value = 163
$ws.Range("J4").Value = 30
$ws.Range("J6").Value = 123
$ws.Range("J7").Value = 430

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 71
$ws.Range("J3").Value = 76
$ws.Range("J7").Value = 359

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J4").Value = 17
$ws.Range("J6").Value = 34
$ws.Range("J7").Value = 128

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J4").Value = 12
$ws.Range("J7").Value = 96

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J3").Value = 44
$ws.Range("J5").Value = 5
$ws.Range("J7").Value = 104

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("J6").Value = 9
$ws.Range("J7").Value = 14

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 80
$ws.Range("J3").Value = 104
$ws.Range("J7").Value = 281

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J2").Value = 61
$ws.Range("J3").Value = 57

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J3").Value = 35
$ws.Range("J7").Value = 137

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J2").Value = 18
$ws.Range("J6").Value = 42
$ws.Range("J7").Value = 77

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J2").Value = 45
$ws.Range("J4").Value = 8
$ws.Range("J6").Value = 48
$ws.Range("J7").Value = 130

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J2").Value = 29
$ws.Range("J7").Value = 80

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J2").Value = 14
$ws.Range("J7").Value = 57

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("J2").Value = 7
$ws.Range("J7").Value = 29

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J2").Value = 21
$ws.Range("J7").Value = 93

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J3").Value = 187
$ws.Range("J6").Value = 166
$ws.Range("J7").Value = 577

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("J6").Value = 14
$ws.Range("J7").Value = 23

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J3").Value = 11
$ws.Range("J4").Value = 26
$ws.Range("J7").Value = 55

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J3").Value = 26
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 101

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("J3").Value = 10
$ws.Range("J7").Value = 43

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J6").Value = 51
$ws.Range("J7").Value = 84

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J6").Value = 53
$ws.Range("J7").Value = 88

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("J6").Value = 7
$ws.Range("J7").Value = 34

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J3").Value = 80
$ws.Range("J6").Value = 97
$ws.Range("J7").Value = 275
